$p = $ppt.ActivePresentation

# Delete slide 4 (the "#103" / newSys1 Func1/Func2 example slide).
$p.Slides.Item(4).Delete()
